$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 999.6
$ws.Cells.Item(28, 9).Value = 999.6
$ws.Cells.Item(28, 11).Value = 999.6
$ws.Cells.Item(28, 13).Value = -514.6

$ws.Cells.Item(53, 8).Value = 315.5
$ws.Cells.Item(53, 10).Value = 449.2
$ws.Cells.Item(53, 12).Value = 449.2
$ws.Cells.Item(53, 14).Value = -1723.2

$ws.Cells.Item(106, 8).Value = 35998
$ws.Cells.Item(106, 9).Value = 35998
$ws.Cells.Item(106, 11).Value = 35998
$ws.Cells.Item(106, 13).Value = -35367

$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()

$ws.Cells.Item(111, 8).Value = 1999
$ws.Cells.Item(111, 10).Value = 1999
$ws.Cells.Item(111, 12).Value = 5997
$ws.Cells.Item(111, 14).Value = -12131

$ws.Cells.Item(125, 8).Value = 166672830
$ws.Cells.Item(125, 9).Value = 142861650
$ws.Cells.Item(125, 11).Value = 1285754850
$ws.Cells.Item(125, 13).Value = -1285752390

$ws.Cells.Item(135, 8).Value = 509.07144
$ws.Cells.Item(135, 9).Value = 394.3846
$ws.Cells.Item(135, 11).Value = 3549.4614
$ws.Cells.Item(135, 13).Value = -1014.4614

$ws.Cells.Item(137, 8).Value = 1275.6
$ws.Cells.Item(137, 9).Value = 1186
$ws.Cells.Item(137, 10).Value = 1410
$ws.Cells.Item(137, 11).Value = 3558
$ws.Cells.Item(137, 12).Value = 4230
$ws.Cells.Item(137, 13).Value = -1008
$ws.Cells.Item(137, 14).Value = -9330

$ws.Cells.Item(138, 8).Value = 4933.5654
$ws.Cells.Item(138, 9).Value = 2984.2856
$ws.Cells.Item(138, 10).Value = 5786.375
$ws.Cells.Item(138, 11).Value = 8952.856800000001
$ws.Cells.Item(138, 12).Value = 17359.125
$ws.Cells.Item(138, 13).Value = -3812.856800000001
$ws.Cells.Item(138, 14).Value = -27639.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3100
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 14).Value = -4424

$ws.Cells.Item(136, 8).Value = 3100
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(127, 8).Value = 68999
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 68999
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 68999
$ws.Cells.Item(127, 14).Value = -78919
$ws.Cells.Item(127, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 66223
$ws.Cells.Item(132, 10).Value = 66223
$ws.Cells.Item(132, 12).Value = 66223
$ws.Cells.Item(132, 14).Value = -76343

$ws.Cells.Item(134, 8).Value = 2673.25
$ws.Cells.Item(134, 9).Value = 2673.25
$ws.Cells.Item(134, 11).Value = 8019.75
$ws.Cells.Item(134, 13).Value = -5484.75

$ws.Cells.Item(141, 8).Value = 86332.336
$ws.Cells.Item(141, 9).Value = 90000
$ws.Cells.Item(141, 10).Value = 84498.5
$ws.Cells.Item(141, 11).Value = 90000
$ws.Cells.Item(141, 12).Value = 84498.5
$ws.Cells.Item(141, 13).Value = -84820
$ws.Cells.Item(141, 14).Value = -94858.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2033.4286
$ws.Cells.Item(58, 10).Value = 2111
$ws.Cells.Item(58, 12).Value = 2111
$ws.Cells.Item(58, 14).Value = -2517

$ws.Cells.Item(86, 8).Value = 19999.666
$ws.Cells.Item(86, 9).Value = 19999.5
$ws.Cells.Item(86, 10).Value = 20000
$ws.Cells.Item(86, 11).Value = 19999.5
$ws.Cells.Item(86, 12).Value = 20000
$ws.Cells.Item(86, 13).Value = -18876.5
$ws.Cells.Item(86, 14).Value = -22246

$ws.Cells.Item(89, 8).Value = 19999.666
$ws.Cells.Item(89, 9).Value = 19999.5
$ws.Cells.Item(89, 10).Value = 20000
$ws.Cells.Item(89, 11).Value = 99997.5
$ws.Cells.Item(89, 12).Value = 100000
$ws.Cells.Item(89, 13).Value = -94381.5
$ws.Cells.Item(89, 14).Value = -111232

$ws.Cells.Item(94, 8).Value = 5458.25
$ws.Cells.Item(94, 9).Value = 6111.75
$ws.Cells.Item(94, 10).Value = 4804.75
$ws.Cells.Item(94, 11).Value = 6111.75
$ws.Cells.Item(94, 12).Value = 4804.75
$ws.Cells.Item(94, 13).Value = -5660.75
$ws.Cells.Item(94, 14).Value = -5706.75

$ws.Cells.Item(99, 8).Value = 5415.5
$ws.Cells.Item(99, 9).Value = 2666.6
$ws.Cells.Item(99, 10).Value = 9997
$ws.Cells.Item(99, 11).Value = 2666.6
$ws.Cells.Item(99, 12).Value = 9997
$ws.Cells.Item(99, 13).Value = -1168.6
$ws.Cells.Item(99, 14).Value = -12993

$ws.Cells.Item(126, 8).Value = 5415.5
$ws.Cells.Item(126, 9).Value = 2666.6
$ws.Cells.Item(126, 10).Value = 9997
$ws.Cells.Item(126, 11).Value = 7999.799999999999
$ws.Cells.Item(126, 12).Value = 29991
$ws.Cells.Item(126, 13).Value = -5529.799999999999
$ws.Cells.Item(126, 14).Value = -34931

$ws.Cells.Item(133, 8).Value = 42647.5
$ws.Cells.Item(133, 10).Value = 59999
$ws.Cells.Item(133, 12).Value = 59999
$ws.Cells.Item(133, 14).Value = -65059

$ws.Cells.Item(136, 8).Value = 2033.4286
$ws.Cells.Item(136, 10).Value = 2111
$ws.Cells.Item(136, 12).Value = 6333
$ws.Cells.Item(136, 14).Value = -11433

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 22001292
$ws.Cells.Item(4, 9).Value = 22001292
$ws.Cells.Item(4, 11).Value = 66003876
$ws.Cells.Item(4, 13).Value = -66003764

$ws.Cells.Item(22, 8).Value = 2450
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 2450
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 7350
$ws.Cells.Item(22, 14).Value = -7688
$ws.Cells.Item(22, 13).ClearContents()

$ws.Cells.Item(27, 8).Value = 2450
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 2450
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 7350
$ws.Cells.Item(27, 14).Value = -7554
$ws.Cells.Item(27, 13).ClearContents()

$ws.Cells.Item(104, 8).Value = 6998
$ws.Cells.Item(104, 10).Value = 6998
$ws.Cells.Item(104, 12).Value = 20994
$ws.Cells.Item(104, 14).Value = -26236

$ws.Cells.Item(129, 8).Value = 4575
$ws.Cells.Item(129, 10).Value = 6333.3335
$ws.Cells.Item(129, 12).Value = 19000.0005
$ws.Cells.Item(129, 14).Value = -29000.0005

$ws.Cells.Item(131, 8).Value = 2570.6667
$ws.Cells.Item(131, 9).Value = 1349.5
$ws.Cells.Item(131, 10).Value = 2814.9
$ws.Cells.Item(131, 11).Value = 4048.5
$ws.Cells.Item(131, 12).Value = 8444.700000000001
$ws.Cells.Item(131, 13).Value = 991.5
$ws.Cells.Item(131, 14).Value = -18524.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 53.24
$ws.Cells.Item(2, 9).Value = 62.1
$ws.Cells.Item(2, 11).Value = 62.1
$ws.Cells.Item(2, 13).Value = 50.9

$ws.Cells.Item(111, 8).Value = 44444
$ws.Cells.Item(111, 10).Value = 44444
$ws.Cells.Item(111, 12).Value = 44444
$ws.Cells.Item(111, 14).Value = -50578

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 533.9
$ws.Cells.Item(16, 9).Value = 504.875
$ws.Cells.Item(16, 11).Value = 504.875
$ws.Cells.Item(16, 13).Value = -334.875

$ws.Cells.Item(132, 8).Value = 13499.8
$ws.Cells.Item(132, 9).Value = 19164.334
$ws.Cells.Item(132, 10).Value = 5003
$ws.Cells.Item(132, 11).Value = 57493.00199999999
$ws.Cells.Item(132, 12).Value = 15009
$ws.Cells.Item(132, 13).Value = -54963.00199999999
$ws.Cells.Item(132, 14).Value = -20069

$ws.Cells.Item(136, 8).Value = 2116.9092
$ws.Cells.Item(136, 9).Value = 1476.3334
$ws.Cells.Item(136, 11).Value = 4429.0002
$ws.Cells.Item(136, 13).Value = -1879.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2057.3333
$ws.Cells.Item(122, 9).Value = 2057.3333
$ws.Cells.Item(122, 11).Value = 6171.999899999999
$ws.Cells.Item(122, 13).Value = -3721.999899999999

$ws.Cells.Item(138, 8).Value = 79999
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 79999
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 79999
$ws.Cells.Item(138, 14).Value = -90279
$ws.Cells.Item(138, 13).ClearContents()
